# Update res_bus/vm_pu.xlsx: bus 0 (slack) voltage setpoint changed from
# 1.05 pu to 1.02 pu ("case with 380 kV done"), which re-ran the power
# flow and changed every bus voltage magnitude result (rows 2-25,
# columns B-N, column H is an out-of-service bus and has no value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.054759570309093
$ws.Range("D2").Value = 1.061497355633386
$ws.Range("E2").Value = 1.06116292618961
$ws.Range("F2").Value = 1.071593373754983
$ws.Range("I2").Value = 1.041099111415091
$ws.Range("J2").Value = 1.059769982955431
$ws.Range("K2").Value = 1.064221255427643
$ws.Range("L2").Value = 1.063887734543676
$ws.Range("M2").Value = 1.074290133938534
$ws.Range("N2").Value = 1.023674500735351

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05596579106465
$ws.Range("D3").Value = 1.062590946997825
$ws.Range("E3").Value = 1.06224397925449
$ws.Range("F3").Value = 1.072781560225625
$ws.Range("I3").Value = 1.041327091674347
$ws.Range("J3").Value = 1.060626370296301
$ws.Range("K3").Value = 1.065128776527847
$ws.Range("L3").Value = 1.064782682965483
$ws.Range("M3").Value = 1.075293986087462
$ws.Range("N3").Value = 1.023961279460058

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056745720884544
$ws.Range("D4").Value = 1.063298368074
$ws.Range("E4").Value = 1.062943373913508
$ws.Range("F4").Value = 1.073550423855197
$ws.Range("I4").Value = 1.041472293771668
$ws.Range("J4").Value = 1.061179397103499
$ws.Range("K4").Value = 1.065715203217502
$ws.Range("L4").Value = 1.065361058248249
$ws.Range("M4").Value = 1.075942994785526
$ws.Range("N4").Value = 1.024146398098472

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.057073467304363
$ws.Range("D5").Value = 1.06359571961203
$ws.Range("E5").Value = 1.063237371850789
$ws.Range("F5").Value = 1.073873661821496
$ws.Range("I5").Value = 1.04153278222928
$ws.Range("J5").Value = 1.06141162396742
$ws.Range("K5").Value = 1.065961546316627
$ws.Range("L5").Value = 1.065604036605413
$ws.Range("M5").Value = 1.076215707146145
$ws.Range("N5").Value = 1.024224115344252

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.057128489441145
$ws.Range("D6").Value = 1.063645643439866
$ws.Range("E6").Value = 1.063286733806972
$ws.Range("F6").Value = 1.073927935381587
$ws.Range("I6").Value = 1.041542905992046
$ws.Range("J6").Value = 1.061450600329355
$ws.Range("K6").Value = 1.066002897250512
$ws.Range("L6").Value = 1.065644823760344
$ws.Range("M6").Value = 1.076261489103325
$ws.Range("N6").Value = 1.024237158165654

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056750100780787
$ws.Range("D7").Value = 1.063302341488357
$ws.Range("E7").Value = 1.062947302435593
$ws.Range("F7").Value = 1.073554742945126
$ws.Range("I7").Value = 1.041473104200198
$ws.Range("J7").Value = 1.061182501172822
$ws.Range("K7").Value = 1.065718495614232
$ws.Range("L7").Value = 1.065364305606308
$ws.Range("M7").Value = 1.075946639294598
$ws.Range("N7").Value = 1.024147436978816

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055167338507349
$ws.Range("D8").Value = 1.061866983117207
$ws.Range("E8").Value = 1.061528298313268
$ws.Range("F8").Value = 1.071994921904546
$ws.Range("I8").Value = 1.041176638060264
$ws.Range("J8").Value = 1.060059634059579
$ws.Range("K8").Value = 1.064528122638999
$ws.Range("L8").Value = 1.064190335357513
$ws.Range("M8").Value = 1.074629505093601
$ws.Range("N8").Value = 1.023771511513977

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052373821406162
$ws.Range("D9").Value = 1.05933607935808
$ws.Range("E9").Value = 1.05902687918682
$ws.Range("F9").Value = 1.069246478248357
$ws.Range("I9").Value = 1.040636483392507
$ws.Range("J9").Value = 1.058072434180156
$ws.Range("K9").Value = 1.062424364971509
$ws.Range("L9").Value = 1.062116125775198
$ws.Range("M9").Value = 1.07230428414402
$ws.Range("N9").Value = 1.023105653338298

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050508345661071
$ws.Range("D10").Value = 1.057647656799603
$ws.Range("E10").Value = 1.057358564398362
$ws.Range("F10").Value = 1.067414215855303
$ws.Range("I10").Value = 1.040264444126723
$ws.Range("J10").Value = 1.056741816857835
$ws.Range("K10").Value = 1.061017652934637
$ws.Range("L10").Value = 1.060729547727576
$ws.Range("M10").Value = 1.070751206170105
$ws.Range("N10").Value = 1.022659425791189

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049699805067866
$ws.Range("D11").Value = 1.056916262048044
$ws.Range("E11").Value = 1.056635984348711
$ws.Range("F11").Value = 1.066620816229911
$ws.Range("I11").Value = 1.040100513685114
$ws.Range("J11").Value = 1.056164250856012
$ws.Range("K11").Value = 1.060407518951431
$ws.Range("L11").Value = 1.060128235700273
$ws.Range("M11").Value = 1.070077995598391
$ws.Range("N11").Value = 1.022465649560187

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049399357916164
$ws.Range("D12").Value = 1.056644543439318
$ws.Range("E12").Value = 1.056367556192924
$ws.Range("F12").Value = 1.066326107942417
$ws.Range("I12").Value = 1.040039196400047
$ws.Range("J12").Value = 1.055949505457707
$ws.Range("K12").Value = 1.060180733669638
$ws.Range("L12").Value = 1.059904742926081
$ws.Range("M12").Value = 1.069827825950811
$ws.Range("N12").Value = 1.02239358840361

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049463810281686
$ws.Range("D13").Value = 1.056702830110035
$ws.Range("E13").Value = 1.056425136312495
$ws.Range("F13").Value = 1.066389324093552
$ws.Range("I13").Value = 1.040052368458069
$ws.Range("J13").Value = 1.055995578674878
$ws.Range("K13").Value = 1.060229386889159
$ws.Range("L13").Value = 1.059952689183059
$ws.Range("M13").Value = 1.069881493157604
$ws.Range("N13").Value = 1.022409049575455

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049674972457641
$ws.Range("D14").Value = 1.05689380265067
$ws.Range("E14").Value = 1.056613796608692
$ws.Range("F14").Value = 1.066596455649069
$ws.Range("I14").Value = 1.040095453875716
$ws.Range("J14").Value = 1.05614650426023
$ws.Range("K14").Value = 1.06038877596595
$ws.Range("L14").Value = 1.060109764555177
$ws.Range("M14").Value = 1.070057318734792
$ws.Range("N14").Value = 1.022459694679097

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049805060579638
$ws.Range("D15").Value = 1.057011460991963
$ws.Range("E15").Value = 1.056730032459221
$ws.Range("F15").Value = 1.066724075611887
$ws.Range("I15").Value = 1.040121943736183
$ws.Range("J15").Value = 1.056239466405771
$ws.Range("K15").Value = 1.060486960342509
$ws.Range("L15").Value = 1.060206525455777
$ws.Range("M15").Value = 1.070165636153026
$ws.Range("N15").Value = 1.022490887651553

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050561989256672
$ws.Range("D16").Value = 1.057696190685836
$ws.Range("E16").Value = 1.057406515561141
$ws.Range("F16").Value = 1.067466870657577
$ws.Range("I16").Value = 1.040275263895412
$ws.Range("J16").Value = 1.056780118401885
$ws.Range("K16").Value = 1.061058123946267
$ws.Range("L16").Value = 1.060769435445919
$ws.Range("M16").Value = 1.070795869696037
$ws.Range("N16").Value = 1.022672274318033

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051036580829187
$ws.Range("D17").Value = 1.058125623082741
$ws.Range("E17").Value = 1.057830804026282
$ws.Range("F17").Value = 1.067932800408668
$ws.Range("I17").Value = 1.040370678369578
$ws.Range("J17").Value = 1.057118879370274
$ws.Range("K17").Value = 1.061416126182552
$ws.Range("L17").Value = 1.061122288231576
$ws.Range("M17").Value = 1.071191005408387
$ws.Range("N17").Value = 1.022785904047876

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051313327231338
$ws.Range("D18").Value = 1.058376075213979
$ws.Range("E18").Value = 1.058078266195709
$ws.Range("F18").Value = 1.068204567734325
$ws.Range("I18").Value = 1.040426058504506
$ws.Range("J18").Value = 1.057316337932712
$ws.Range("K18").Value = 1.061624844663725
$ws.Range("L18").Value = 1.061328013137984
$ws.Range("M18").Value = 1.071421412250857
$ws.Range("N18").Value = 1.022852128654733

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051407677981519
$ws.Range("D19").Value = 1.058461468124484
$ws.Range("E19").Value = 1.058162641369574
$ws.Range("F19").Value = 1.068297233144394
$ws.Range("I19").Value = 1.040444895315223
$ws.Range("J19").Value = 1.057383643344986
$ws.Range("K19").Value = 1.061695995661564
$ws.Range("L19").Value = 1.061398145101175
$ws.Range("M19").Value = 1.071499963347158
$ws.Range("N19").Value = 1.022874700438718

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050985669403775
$ws.Range("D20").Value = 1.058079552012273
$ws.Range("E20").Value = 1.057785283774093
$ws.Range("F20").Value = 1.067882810735966
$ws.Range("I20").Value = 1.040360469592821
$ws.Range("J20").Value = 1.057082547496514
$ws.Range("K20").Value = 1.061377726095191
$ws.Range("L20").Value = 1.061084439595533
$ws.Range("M20").Value = 1.071148618248305
$ws.Range("N20").Value = 1.022773718209989

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049612793769049
$ws.Range("D21").Value = 1.056837567279892
$ws.Range("E21").Value = 1.056558241686303
$ws.Range("F21").Value = 1.066535460701723
$ws.Range("I21").Value = 1.040082778062287
$ws.Range("J21").Value = 1.056102066291062
$ws.Range("K21").Value = 1.060341844127373
$ws.Range("L21").Value = 1.060063513612845
$ws.Range("M21").Value = 1.070005545512748
$ws.Range("N21").Value = 1.022444783281344

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048748921659031
$ws.Range("D22").Value = 1.056056415116776
$ws.Range("E22").Value = 1.055786578783447
$ws.Range("F22").Value = 1.065688302318687
$ws.Range("I22").Value = 1.039905716242745
$ws.Range("J22").Value = 1.055484373511818
$ws.Range("K22").Value = 1.059689650293748
$ws.Range("L22").Value = 1.059420813713696
$ws.Range("M22").Value = 1.069286217398561
$ws.Range("N22").Value = 1.022237482430466

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049206942885809
$ws.Range("D23").Value = 1.056470544466968
$ws.Range("E23").Value = 1.056195668742539
$ws.Range("F23").Value = 1.066137400019959
$ws.Range("I23").Value = 1.039999813926241
$ws.Range("J23").Value = 1.055811940742394
$ws.Range("K23").Value = 1.060035475755657
$ws.Range("L23").Value = 1.059761597666608
$ws.Range("M23").Value = 1.069667607327358
$ws.Range("N23").Value = 1.022347422770975

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051008674316709
$ws.Range("D24").Value = 1.058100369634138
$ws.Range("E24").Value = 1.057805852472737
$ws.Range("F24").Value = 1.067905398922372
$ws.Range("I24").Value = 1.040365083344781
$ws.Range("J24").Value = 1.057098964723513
$ws.Range("K24").Value = 1.061395077745153
$ws.Range("L24").Value = 1.061101542037445
$ws.Range("M24").Value = 1.071167771395109
$ws.Range("N24").Value = 1.02277922463153

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.053096556005975
$ws.Range("D25").Value = 1.059990577121224
$ws.Range("E25").Value = 1.059673675158606
$ws.Range("F25").Value = 1.069957005419397
$ws.Range("I25").Value = 1.040778229628665
$ws.Range("J25").Value = 1.058587193950262
$ws.Range("K25").Value = 1.062968973827273
$ws.Range("L25").Value = 1.062653019174368
$ws.Range("M25").Value = 1.072905921224816
$ws.Range("N25").Value = 1.023278201708266
